# Daily attendance processing - swap order of first two "Recorded By" entries
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -ge 2) {
            $newParts = @($parts[1], $parts[0])
            if ($parts.Count -gt 2) {
                $newParts += $parts[2..($parts.Count - 1)]
            }
            $newVal = $newParts -join ", "
            $cell.Value = $newVal
        }
    }
}
